$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Recon_Config (sheet1): insert a new leading column "Use Case" /
#    "Legacy_vs_New" so every existing column shifts one to the right.
# ---------------------------------------------------------------------------
$wsRecon = $wb.Worksheets.Item("Recon_Config")

$wsRecon.Columns.Item(1).Insert()
$wsRecon.Range("A1").Value = "Use Case"
$wsRecon.Range("A2").Value = "Legacy_vs_New"

# Match the formatting of the (now shifted) neighboring header/data cells.
$wsRecon.Range("B1").Copy()
$wsRecon.Range("A1").PasteSpecial(-4122)
$wsRecon.Range("B2").Copy()
$wsRecon.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsRecon.Columns.Item(1).ColumnWidth = 13.25

# ---------------------------------------------------------------------------
# 2) Insert a brand-new worksheet "Legacy_vs_New" between Recon_Config and
#    Config_Usage, carrying the Source/Target column mapping table.
# ---------------------------------------------------------------------------
$wsLegacy = $wb.Worksheets.Add($null, $wsRecon)
$wsLegacy.Name = "Legacy_vs_New"

$wsLegacy.Range("A1").Value = "Source_Column"
$wsLegacy.Range("B1").Value = "Target_Column"
$wsLegacy.Range("A2").Value = "Store_Number"
$wsLegacy.Range("B2").Value = "Store_ID"
$wsLegacy.Range("A3").Value = "XYZ"
$wsLegacy.Range("B3").Value = "ABC"

# Give the table the same bordered-cell look used elsewhere in the workbook.
$wsRecon.Range("B2").Copy()
$wsLegacy.Range("A1:B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsLegacy.Range("A1:B3").EntireColumn.ColumnWidth = 12.92
$wsLegacy.Range("D8").Select()

# ---------------------------------------------------------------------------
# 3) Config_Usage: only the remembered selection changed.
# ---------------------------------------------------------------------------
$wsUsage = $wb.Worksheets.Item("Config_Usage")
$wsUsage.Range("A32").Select()

# ---------------------------------------------------------------------------
# 4) Restore Recon_Config as the active sheet/selection.
# ---------------------------------------------------------------------------
$wsRecon.Activate()
$wsRecon.Range("U7").Select()
